$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# 1. Insert a new column before column D ("PopulationId" stays C, everything from
#    old D ("ModelParameterSheets") onward shifts one column to the right).
$ws.Columns("D").Insert()

# 2. Add the new header in D1, matching the style (bold) of the other header cells.
$ws.Range("D1").Value = "ReadPopulationFromCSV"
$ws.Range("D1").Font.Bold = $true

# 3. Row 2 (TestScenario) previously had SteadyState = FALSE; that value shifted
#    from H2 to I2 with the column insert. The new layout has this cell blank.
$ws.Range("I2").ClearContents()

# 4. Row 4 (PopulationScenario) gets the new ReadPopulationFromCSV value = FALSE.
$ws.Range("D4").Value = $false

# 5. Add new row 5 for the PopulationScenarioFromCSV scenario.
$ws.Range("A5").Value = "PopulationScenarioFromCSV"
$ws.Range("B5").Value = "Indiv"
$ws.Range("C5").Value = "TestPopulation"
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "Global"
$ws.Range("F5").Value = "Aciclovir_iv_250mg"
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = "h"
$ws.Range("I5").Value = $false
$ws.Range("L5").Value = "Aciclovir.pkml"

# 6. Autofit the columns to match the widths Excel would compute for the new content.
$ws.Columns("A:M").AutoFit()

# 7. Update the selection to match the final state.
$ws.Range("G13").Select()
